$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.510.80'
$ws.Range("E2").Value = '  -2.00%  '

# Row 3
$ws.Range("D3").Value = '3.475.12'
$ws.Range("E3").Value = '  -4.51%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.98'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -4.50%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '192.15'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -4.27%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.609'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -2.98%  '

# Row 8
$ws.Range("D8").Value = '3.461.95'
$ws.Range("E8").Value = '  -4.57%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.08%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.206'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -6.21%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.618'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -4.41%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '51.44'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -4.67%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000286'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -6.83%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.12'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -4.84%  '

# Row 15
$ws.Range("D15").Value = '4.009.92'
$ws.Range("E15").Value = '  -4.71%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '643.23'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.94%  '

# Row 17
$ws.Range("D17").Value = '69.221.83'
$ws.Range("E17").Value = '  -2.39%  '

# Row 18
$ws.Range("D18").Value = '3.446.18'
$ws.Range("E18").Value = '  -4.81%  '

# Row 19
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.28'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -5.72%  '

# Row 20
$ws.Range("B20").Value = 'TRON'
$ws.Range("C20").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.121'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.95%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '18.12'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -5.07%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.942'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -5.89%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.84'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.45%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.30'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.40%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '99.07'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -4.95%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.27'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -7.78%  '

# Row 27
$ws.Range("E27").Value = '  -5.11%  '

# Row 28
$ws.Range("E28").Value = '  -6.86%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.28'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -4.96%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.34'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -4.36%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.32'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -9.85%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.72'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -6.76%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.56'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -5.88%  '

# Row 34
$ws.Range("E34").Value = '  -6.26%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '61.03'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -3.75%  '

# Row 36
$ws.Range("B36").Value = 'Maker'
$ws.Range("C36").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D36").Value = '3.691.48'
$ws.Range("E36").Value = '  -7.55%  '

# Row 37
$ws.Range("B37").Value = 'Dai'
$ws.Range("C37").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.05%  '

# Row 38
$ws.Range("D38").Value = '0.0₃0794'
$ws.Range("E38").Value = '  -10.00%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '504.96'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.25%  '

# Row 40
$ws.Range("E40").Value = '  -3.71%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.49'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.57%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.370'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -5.10%  '

# Row 43
$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.133'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.55%  '

# Row 44
$ws.Range("B44").Value = 'CoreDAO'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.51'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +69.46%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '34.31'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -6.71%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0442'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -4.63%  '

# Row 47
$ws.Range("E47").Value = '  -4.88%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.82'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -3.58%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.134'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -4.98%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.998'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.32%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.11'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -6.13%  '
